# Generate Report for Handback
# Update the timestamp strings recorded for the handoff/handback of the
# e07622b0-cb26-40e5-98df-890a223460f3 file across the Overview, zh-cn and
# de-de worksheets to reflect the latest report generation.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" for e07622b0-...md
$wsOverview.Range("G4").Value = "2016-09-01 16:54:16"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H4").Value = "2016-09-01 16:54:10"
$wsZhCn.Range("K4").Value = "2016-09-01 16:54:38"

# de-de sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsDeDe.Range("H4").Value = "2016-09-01 16:54:16"
$wsDeDe.Range("K4").Value = "2016-09-01 16:54:47"
